$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text interpretation for numeric-looking price cells so they
# keep their original display (e.g. trailing zeros, thousand-dot grouping)
# instead of being auto-converted to Excel numbers.
$ws.Range("D2").Value = '26.998.20'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '1.560.00'
$ws.Range("E3").Value = '  +0.46%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '208.36'
$ws.Range("E5").Value = '  +0.74%  '
$ws.Range("E6").Value = '  +0.41%  '
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.01'
$ws.Range("E8").Value = '  -0.61%  '
$ws.Range("E9").Value = '  +0.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0597'
$ws.Range("E10").Value = '  +1.75%  '
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("D12").Value = '1.781.78'
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("D13").Value = '1.504.77'
$ws.Range("E13").Value = '  -3.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.73'
$ws.Range("E14").Value = '  -0.30%  '
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("D16").Value = '26.985.08'
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.85'
$ws.Range("E17").Value = '  +0.26%  '
$ws.Range("E18").Value = '  +1.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.81'
$ws.Range("E19").Value = '  -0.78%  '
$ws.Range("E20").Value = '  +1.14%  '
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.14'
$ws.Range("E22").Value = '  +1.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.20'
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.94'
$ws.Range("E24").Value = '  -0.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.81'
$ws.Range("E25").Value = '  -0.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.60'
$ws.Range("E26").Value = '  -0.35%  '
$ws.Range("E27").Value = '  +0.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.106'
$ws.Range("E28").Value = '  +1.47%  '
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0474'
$ws.Range("E30").Value = '  +1.63%  '
$ws.Range("E31").Value = '  +3.71%  '
$ws.Range("E32").Value = '  +0.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.17'
$ws.Range("D34").Value = '1.424.49'
$ws.Range("E34").Value = '  +0.17%  '
$ws.Range("B35").Value = 'TrustWalletToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.07'
$ws.Range("E35").Value = '  +10.40%  '
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.60'
$ws.Range("E36").Value = '  +1.58%  '
$ws.Range("E37").Value = '  +2.19%  '
$ws.Range("E38").Value = '  +0.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.533'
$ws.Range("E39").Value = '  +2.39%  '
$ws.Range("E40").Value = '  +2.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.809'
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.01'
$ws.Range("E42").Value = '  +0.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.88%  '
$ws.Range("E44").Value = '  -0.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.59'
$ws.Range("E45").Value = '  +0.45%  '
$ws.Range("E46").Value = '  -0.84%  '
$ws.Range("D47").Value = '1.696.13'
$ws.Range("E47").Value = '  +0.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.76'
$ws.Range("E48").Value = '  -1.08%  '
$ws.Range("D49").Value = '0.0₆0103'
$ws.Range("E49").Value = '  +3.10%  '
$ws.Range("E50").Value = '  -0.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0958'
$ws.Range("E51").Value = '  +0.57%  '
